# AST Interpreter and tests (#50)
#
# Adds a new "Interpreter" compiler stage (and a "Compiler" stage, used
# elsewhere) to the Stages lookup sheet, and a corresponding new error row
# ("Accessed a null variable (unassigned)") to the Errors table.

$wb = $excel.ActiveWorkbook

# --- Stages sheet: register the new "Compiler" / "Interpreter" stages ---
$wsStages = $wb.Worksheets.Item("Stages")
$wsStages.Range("A6").Value = "Compiler"
$wsStages.Range("B6").Value = 68000
$wsStages.Range("A7").Value = "Interpreter"
$wsStages.Range("B7").Value = 69000
# Leave the selection where a user would land after typing the two rows.
$wsStages.Range("B8").Select() | Out-Null

# --- Errors sheet: grow Table1 by one row and fill in the new error ---
$wsErrors = $wb.Worksheets.Item("Errors")
$lo = $wsErrors.ListObjects.Item("Table1")
$lo.Resize($wsErrors.Range("A1:F29"))

$wsErrors.Range("A29").Value = "Error"
$wsErrors.Range("B29").Value = 1
$wsErrors.Range("C29").Value = "Interpreter"
$wsErrors.Range("D29").Value = "Accessed a null variable (unassigned)"
$wsErrors.Range("E29").Formula = '=(_xlfn.XLOOKUP($C29,Stages!$A:$A,Stages!$B:$B)+$B29)'
$wsErrors.Range("F29").Formula = '=LEFT(A29,1)&E29'

# Re-activate the Errors sheet/selection (matches the author's final state).
$wsErrors.Range("D30").Select() | Out-Null
